$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the stray _GoBack bookmark near the top of the document
#    (it gets re-added later, around the new "DX Grid" text).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2. Delete the whole first bullet ("Click on waterfall to set QSO
#    Frequency. ... CTRL key during the move.") and merge its
#    paragraph into the following one.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Click on waterfall to set QSO Frequency*the CTRL key during the move.", $false, $false, $true, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find first bullet text to delete"
}
$rng.Delete()

# Now delete the trailing paragraph mark left behind so the paragraphs merge.
$rng.Collapse(0)
$rng.MoveEnd(1, 1)
if ($rng.Text -eq [char]13) {
    $rng.Delete()
}

Write-Output "done stage 1/2"
